# Add data for 2022-11-09
#
# - Rename the sheet to reflect the new "through" date (2022-11-01)
# - Relabel row 11 from "October (through 10-31)" to plain "October"
#   (its data, the October-only counts, does not change)
# - Row 12 used to be the "Total" row (cumulative through 10-31). It is
#   repurposed as the new "November (through 11-01)" row holding just the
#   new November counts.
# - Add a brand new row 13 as the "Total" row: the old Total-through-October
#   values plus the new November counts.
# - Widen column A slightly (23.71 -> 24.71 characters) to fit better.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wb.Worksheets.Item(1).Name = "Through 2022-11-01"

# Row 11: "October (through 10-31)" -> "October" (values unchanged)
$ws.Cells.Item(11, 1).Value = "October"

# Row 12: was the Total-through-Oct-31 row; becomes the November row.
$ws.Cells.Item(12, 1).Value = "November (through 11-01)"
$ws.Cells.Item(12, 2).Value = 1
$ws.Cells.Item(12, 3).Value = 3
$ws.Cells.Item(12, 4).Value = 5
$ws.Cells.Item(12, 5).Value = 4
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 9
$ws.Cells.Item(12, 8).Value = 5
$ws.Cells.Item(12, 9).Value = 2

# Row 13 (new): the Total row = old Total-through-Oct values + new Nov row.
$ws.Cells.Item(13, 1).Value = "Total"
$ws.Cells.Item(13, 2).Value = 259
$ws.Cells.Item(13, 3).Value = 489
$ws.Cells.Item(13, 4).Value = 715
$ws.Cells.Item(13, 5).Value = 619
$ws.Cells.Item(13, 6).Value = 483
$ws.Cells.Item(13, 7).Value = 1066
$ws.Cells.Item(13, 8).Value = 1446
$ws.Cells.Item(13, 9).Value = 1403

# Match the bold/centered/bordered label style used by the other month
# label cells in column A (row 1 header style == row 11/12 label style).
$ws.Cells.Item(13, 1).Font.Bold = $true
$ws.Cells.Item(13, 1).HorizontalAlignment = -4108
$ws.Cells.Item(13, 1).VerticalAlignment = -4160
$ws.Cells.Item(13, 1).Borders.LineStyle = 1

$ws.Columns.Item(1).ColumnWidth = 24.7109375
